$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, new DAMSLTag (col I), new DialogAct (col J)
$updates = @(
    ,@(10, 'sv', 'Statement-opinion')
    ,@(12, 'sd', 'Statement-non-opinion')
    ,@(38, 'sd', 'Statement-non-opinion')
    ,@(41, 'sv', 'Statement-opinion')
    ,@(68, 'b', 'Acknowledge (Backchannel)')
    ,@(74, 'sd', 'Statement-non-opinion')
    ,@(90, 'sd', 'Statement-non-opinion')
    ,@(102, 'sd', 'Statement-non-opinion')
    ,@(103, 'aa', 'Agree/Accept')
    ,@(110, 'sd', 'Statement-non-opinion')
    ,@(124, 'sv', 'Statement-opinion')
    ,@(126, 'sv', 'Statement-opinion')
    ,@(127, 'sd', 'Statement-non-opinion')
    ,@(128, 'sv', 'Statement-opinion')
    ,@(129, 'aa', 'Agree/Accept')
    ,@(174, 'sd', 'Statement-non-opinion')
    ,@(183, '%', 'Uninterpretable')
    ,@(184, 'b', 'Acknowledge (Backchannel)')
    ,@(187, 'aa', 'Agree/Accept')
    ,@(191, 'sd', 'Statement-non-opinion')
    ,@(199, 'sv', 'Statement-opinion')
    ,@(210, 'sv', 'Statement-opinion')
    ,@(213, 'sv', 'Statement-opinion')
    ,@(221, 'aa', 'Agree/Accept')
    ,@(224, 'aa', 'Agree/Accept')
    ,@(226, '%', 'Uninterpretable')
    ,@(259, 'aa', 'Agree/Accept')
    ,@(276, 'aa', 'Agree/Accept')
    ,@(280, 'sv', 'Statement-opinion')
    ,@(286, 'b', 'Acknowledge (Backchannel)')
    ,@(288, 'b', 'Acknowledge (Backchannel)')
    ,@(290, 'aa', 'Agree/Accept')
    ,@(291, 'b', 'Acknowledge (Backchannel)')
    ,@(293, 'sd', 'Statement-non-opinion')
    ,@(300, 'aa', 'Agree/Accept')
    ,@(328, '%', 'Uninterpretable')
    ,@(335, 'aa', 'Agree/Accept')
    ,@(351, 'sv', 'Statement-opinion')
    ,@(353, 'sv', 'Statement-opinion')
    ,@(355, 'sv', 'Statement-opinion')
    ,@(356, 'sd', 'Statement-non-opinion')
    ,@(365, 'aa', 'Agree/Accept')
    ,@(371, 'sd', 'Statement-non-opinion')
    ,@(385, 'aa', 'Agree/Accept')
    ,@(394, 'sd', 'Statement-non-opinion')
    ,@(397, 'aa', 'Agree/Accept')
    ,@(398, 'sv', 'Statement-opinion')
    ,@(402, 'aa', 'Agree/Accept')
    ,@(410, 'sd', 'Statement-non-opinion')
    ,@(412, 'sv', 'Statement-opinion')
    ,@(417, 'sd', 'Statement-non-opinion')
    ,@(420, 'sv', 'Statement-opinion')
    ,@(430, 'sd', 'Statement-non-opinion')
    ,@(431, 'b', 'Acknowledge (Backchannel)')
    ,@(432, 'aa', 'Agree/Accept')
    ,@(435, 'b', 'Acknowledge (Backchannel)')
    ,@(436, 'sd', 'Statement-non-opinion')
    ,@(439, 'b', 'Acknowledge (Backchannel)')
    ,@(442, 'aa', 'Agree/Accept')
    ,@(443, 'b', 'Acknowledge (Backchannel)')
    ,@(445, 'b', 'Acknowledge (Backchannel)')
    ,@(446, 'aa', 'Agree/Accept')
    ,@(448, 'sd', 'Statement-non-opinion')
    ,@(461, 'aa', 'Agree/Accept')
    ,@(489, 'sd', 'Statement-non-opinion')
    ,@(490, 'sd', 'Statement-non-opinion')
    ,@(508, 'sv', 'Statement-opinion')
    ,@(511, 'ba', 'Appreciation')
    ,@(542, 'sv', 'Statement-opinion')
    ,@(543, 'aa', 'Agree/Accept')
    ,@(571, 'aa', 'Agree/Accept')
    ,@(577, 'sd', 'Statement-non-opinion')
    ,@(589, 'aa', 'Agree/Accept')
    ,@(594, 'ba', 'Appreciation')
    ,@(598, 'aa', 'Agree/Accept')
    ,@(599, 'aa', 'Agree/Accept')
    ,@(600, 'aa', 'Agree/Accept')
    ,@(603, 'sd', 'Statement-non-opinion')
    ,@(615, 'aa', 'Agree/Accept')
    ,@(616, 'sd', 'Statement-non-opinion')
    ,@(621, 'sv', 'Statement-opinion')
    ,@(627, 'aa', 'Agree/Accept')
    ,@(636, 'sd', 'Statement-non-opinion')
    ,@(640, '%', 'Uninterpretable')
    ,@(653, 'b', 'Acknowledge (Backchannel)')
    ,@(655, 'ba', 'Appreciation')
    ,@(658, 'sd', 'Statement-non-opinion')
)

foreach ($u in $updates) {
    $row = $u[0]
    $iVal = $u[1]
    $jVal = $u[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}

Write-Output "Updated $($updates.Count) rows."